$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 300
$ws.Range("I11").Value = 300
$ws.Range("K11").Value = 300
$ws.Range("M11").Value = -160
$ws.Range("H17").Value = 4000
$ws.Range("J17").Value = 4000
$ws.Range("L17").Value = 12000
$ws.Range("N17").Value = -12336
$ws.Range("H40").Value = 4549.8335
$ws.Range("H58").Value = 5737.125
$ws.Range("I58").Value = 3779.4
$ws.Range("J58").Value = 9000
$ws.Range("K58").Value = 11338.2
$ws.Range("L58").Value = 27000
$ws.Range("M58").Value = -11188.2
$ws.Range("N58").Value = -27300
$ws.Range("H107").Value = 2038.2307
$ws.Range("I107").Value = 2038.2307
$ws.Range("K107").Value = 2038.2307
$ws.Range("M107").Value = -118.2307000000001
$ws.Range("H125").Value = 715
$ws.Range("J125").Value = 715
$ws.Range("L125").Value = 6435
$ws.Range("N125").Value = -11355
$ws.Range("H129").Value = 2052.875
$ws.Range("J129").Value = 2086.1724
$ws.Range("L129").Value = 6258.5172
$ws.Range("N129").Value = -16258.5172
$ws.Range("H138").Value = 9678.333000000001
$ws.Range("J138").Value = 9678.333000000001
$ws.Range("L138").Value = 29034.999
$ws.Range("N138").Value = -39314.999
$ws.Range("H141").Value = 1347.2858
$ws.Range("I141").Value = 1347.2858
$ws.Range("K141").Value = 4041.8574
$ws.Range("M141").Value = 1138.1426

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11982.4
$ws.Range("I32").Value = 10302.71
$ws.Range("J32").Value = 25000
$ws.Range("K32").Value = 10302.71
$ws.Range("L32").Value = 25000
$ws.Range("M32").Value = -10015.71
$ws.Range("N32").Value = -25574
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("H63").Value = 4487.5
$ws.Range("J63").Value = 7614
$ws.Range("L63").Value = 7614
$ws.Range("N63").Value = -8986
$ws.Range("H66").Value = 4487.5
$ws.Range("J66").Value = 7614
$ws.Range("L66").Value = 38070
$ws.Range("N66").Value = -44934
$ws.Range("H97").Value = 4019.75
$ws.Range("I97").Value = 722
$ws.Range("J97").Value = 7317.5
$ws.Range("K97").Value = 722
$ws.Range("L97").Value = 7317.5
$ws.Range("M97").Value = -226
$ws.Range("N97").Value = -8309.5
$ws.Range("H132").Value = 3458.2
$ws.Range("I132").Value = 2431
$ws.Range("K132").Value = 7293
$ws.Range("M132").Value = -4763
$ws.Range("H133").Value = 50000
$ws.Range("J133").Value = 50000
$ws.Range("L133").Value = 50000
$ws.Range("N133").Value = -55060
$ws.Range("H135").Value = 93472.336
$ws.Range("J135").Value = 93472.336
$ws.Range("L135").Value = 93472.336
$ws.Range("N135").Value = -103612.336
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("M136").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("N138").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H9").Value = 68750
$ws.Range("J9").Value = 68750
$ws.Range("L9").Value = 68750
$ws.Range("N9").Value = -69086
$ws.Range("H99").Value = 18876.5
$ws.Range("I99").Value = 22151.8
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 22151.8
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = -20653.8
$ws.Range("N99").Value = -5496
$ws.Range("H105").Value = 2179.6667
$ws.Range("I105").Value = 2179.6667
$ws.Range("K105").Value = 2179.6667
$ws.Range("M105").Value = -432.6667000000002
$ws.Range("H134").Value = 2309.3333
$ws.Range("I134").Value = 2339
$ws.Range("K134").Value = 7017
$ws.Range("M134").Value = -4482

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2675.4443
$ws.Range("I31").Value = 2097.75
$ws.Range("K31").Value = 2097.75
$ws.Range("M31").Value = -1802.75
$ws.Range("H34").Value = 2675.4443
$ws.Range("I34").Value = 2097.75
$ws.Range("K34").Value = 2097.75
$ws.Range("M34").Value = -1895.75
$ws.Range("H122").Value = 2853.8235
$ws.Range("I122").Value = 2958.5833
$ws.Range("K122").Value = 8875.749899999999
$ws.Range("M122").Value = -6425.749899999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 146.2
$ws.Range("I10").Value = 146.2
$ws.Range("K10").Value = 438.6
$ws.Range("M10").Value = -299.6
$ws.Range("H131").Value = 1730
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 566.625
$ws.Range("I2").Value = 566.625
$ws.Range("K2").Value = 566.625
$ws.Range("M2").Value = -453.625
$ws.Range("H97").Value = 1639.6
$ws.Range("I97").Value = 1639.6
$ws.Range("K97").Value = 1639.6
$ws.Range("M97").Value = -1143.6

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1663.875
$ws.Range("J16").Value = 5001.5
$ws.Range("L16").Value = 5001.5
$ws.Range("N16").Value = -5341.5
$ws.Range("H22").Value = 8315.380999999999
$ws.Range("I22").Value = 9378.362999999999
$ws.Range("J22").Value = 7146.1
$ws.Range("K22").Value = 9378.362999999999
$ws.Range("L22").Value = 7146.1
$ws.Range("M22").Value = -9083.362999999999
$ws.Range("N22").Value = -7736.1
$ws.Range("H27").Value = 8315.380999999999
$ws.Range("I27").Value = 9378.362999999999
$ws.Range("J27").Value = 7146.1
$ws.Range("K27").Value = 9378.362999999999
$ws.Range("L27").Value = 7146.1
$ws.Range("M27").Value = -9271.362999999999
$ws.Range("N27").Value = -7360.1
$ws.Range("H40").Value = 2156.1428
$ws.Range("I40").Value = 1682.3334
$ws.Range("J40").Value = 4999
$ws.Range("K40").Value = 1682.3334
$ws.Range("L40").Value = 4999
$ws.Range("M40").Value = -1546.3334
$ws.Range("N40").Value = -5271
$ws.Range("H82").Value = 2474.125
$ws.Range("I82").Value = 2860.75
$ws.Range("J82").Value = 2087.5
$ws.Range("K82").Value = 2860.75
$ws.Range("L82").Value = 2087.5
$ws.Range("M82").Value = -2499.75
$ws.Range("N82").Value = -2809.5
$ws.Range("H85").Value = 2474.125
$ws.Range("I85").Value = 2860.75
$ws.Range("J85").Value = 2087.5
$ws.Range("K85").Value = 2860.75
$ws.Range("L85").Value = 2087.5
$ws.Range("M85").Value = -1612.75
$ws.Range("N85").Value = -4583.5
$ws.Range("H93").Value = 2573.75
$ws.Range("I93").Value = 2464.6667
$ws.Range("J93").Value = 2714
$ws.Range("K93").Value = 2464.6667
$ws.Range("L93").Value = 2714
$ws.Range("M93").Value = -1216.6667
$ws.Range("N93").Value = -5210
$ws.Range("H122").Value = 4379.4375
$ws.Range("I122").Value = 3172.9167
$ws.Range("K122").Value = 9518.750100000001
$ws.Range("M122").Value = -7068.750100000001
$ws.Range("H132").Value = 3100.5
$ws.Range("I132").Value = 2972.6191
$ws.Range("J132").Value = 3995.6667
$ws.Range("K132").Value = 8917.8573
$ws.Range("L132").Value = 11987.0001
$ws.Range("M132").Value = -6387.8573
$ws.Range("N132").Value = -17047.0001

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H105").Value = 48000
$ws.Range("J105").Value = 48000
$ws.Range("L105").Value = 48000
$ws.Range("N105").Value = -54988
$ws.Range("H122").Value = 1024.1428
$ws.Range("I122").Value = 1024.1428
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3072.4284
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -622.4284000000002
$ws.Range("H126").Value = 1565.9524
$ws.Range("I126").Value = 996.25
$ws.Range("K126").Value = 2988.75
$ws.Range("M126").Value = -518.75
$ws.Range("H132").Value = 4111.375
$ws.Range("I132").Value = 2128.7
$ws.Range("K132").Value = 6386.099999999999
$ws.Range("M132").Value = -3856.099999999999
$ws.Range("N122").ClearContents()
